$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1. Update the formulas in Sheet1: remove the "*100" factor from the
#    percentage-difference formulas (columns G and I), for rows 8-21, now
#    that the cells will carry a native percentage number format instead.
# ---------------------------------------------------------------------------
for ($r = 8; $r -le 19; $r++) {
    $ws1.Cells.Item($r, 7).Formula = "=IFERROR((F$r-E$r)/E$r,0)"
}
$ws1.Cells.Item(20, 7).Formula = "=IFERROR((F20-E20)/E20,0)"
$ws1.Cells.Item(21, 7).Formula = "=IFERROR((F21-E21)/E21,0)"

for ($r = 8; $r -le 19; $r++) {
    $ws1.Cells.Item($r, 9).Formula = "=IFERROR(H$r/`$E`$21,0)"
}

# ---------------------------------------------------------------------------
# 2. Apply number formats:
#    - E,F,H columns (rows 8-21): accounting-style "#,##0.00"
#    - G,I columns (rows 8-21): percent "0.00%"
# ---------------------------------------------------------------------------
$ws1.Range("E8:F18").NumberFormat = "#,##0.00"
$ws1.Range("H8:H18").NumberFormat = "#,##0.00"
$ws1.Range("G8:G21").NumberFormat = "0.00%"
$ws1.Range("I8:I19").NumberFormat = "0.00%"

$ws1.Range("E19:F19").NumberFormat = "#,##0.00"
$ws1.Range("H19").NumberFormat = "#,##0.00"

$ws1.Range("E20:F20").NumberFormat = "#,##0.00"
$ws1.Range("H20").NumberFormat = "#,##0.00"
$ws1.Range("I20").NumberFormat = "0.00%"

$ws1.Range("E21:F21").NumberFormat = "#,##0.00"
$ws1.Range("H21").NumberFormat = "#,##0.00"
$ws1.Range("I21").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# 3. New Reporters (row 19) - H19 now carries an explicit value instead of
#    being blank.
# ---------------------------------------------------------------------------
$ws1.Range("H19").Value = 616411

# ---------------------------------------------------------------------------
# 4. Sheet view / selection: the sheet had scrolled so A7 was the top-left
#    cell and K20 was selected; now the sheet is scrolled back to the top
#    and H21 is the selected cell.
# ---------------------------------------------------------------------------
$ws1.Range("A1").Select() | Out-Null
$ws1.Range("H21").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Application window geometry (cosmetic - best effort).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Width = 1452
$win.Height = 792
$win.Left = 1434
$win.Top = -6

$wb.Save()
